# Update PRESUPUESTO (column G) values on the "VENTA MENSUAL" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$updates = @{
    3  = 1000
    4  = 15000
    5  = 10000
    6  = 1500
    8  = 500
    9  = 200
    10 = 1500
    11 = 1000
    14 = 3000
    15 = 6000
    16 = 6000
    18 = 3000
    19 = 1000
    20 = 500
    21 = 1000
    28 = 1500
    29 = 6000
    31 = 1500
    32 = 3000
    33 = 300
    35 = 10000
    36 = 8000
    37 = 200
    40 = 3500
    41 = 1000
    42 = 2000
    43 = 200
    45 = 300
    46 = 1500
    47 = 1500
    48 = 1500
    49 = 500
    52 = 1500
    53 = 500
    54 = 1000
    55 = 114500
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
